$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B5").Value = 1128.3373886605052
$ws.Range("C5").Value = 0.45988922593868314
$ws.Range("D5").Value = 10.845094145111563

$ws.Range("B7").Value = 2426.4506570685671
$ws.Range("C7").Value = 0.15232489243450137
$ws.Range("D7").Value = 2.6057672424780725

$ws.Range("B8").Value = 2820.9152543502082
$ws.Range("C8").Value = 0.19932891299178354
$ws.Range("D8").Value = 3.1139391540780683

$ws.Range("B11").Value = 1053.104263588699
$ws.Range("D11").Value = 100.03863232268689

$ws.Range("B12").Value = 1282.4501697944279
$ws.Range("D12").Value = 115.52763125928119

$ws.Range("B14").Value = 1104.5011299338821
$ws.Range("D14").Value = 24.141499809575219

$ws.Range("B17").Value = 5637.8466446024022
$ws.Range("C17").Value = 0.33003484373946029
$ws.Range("D17").Value = 5.7000698096110511

$ws.Range("B18").Value = 5471.8375870387817
$ws.Range("C18").Value = 0.27192622265036603
$ws.Range("D18").Value = 96.544665717026646

$ws.Range("B19").Value = 1488.9441239389866
$ws.Range("D19").Value = 109.74279955399273

$ws.Range("B20").Value = 4192.3992542192427
$ws.Range("C20").Value = 0.19499290911187689
$ws.Range("D20").Value = 6.552522257930308

$ws.Range("B21").Value = 1209.1104117115547
$ws.Range("D21").Value = 183.81922732339979

$ws.Range("B22").Value = 1185.8695325637414
$ws.Range("D22").Value = 94.627258516538916

$ws.Range("B23").Value = 1124.9904217585265
$ws.Range("D23").Value = 109.42911195400781

$ws.Range("B24").Value = 1162.5127380647236
$ws.Range("D24").Value = 110.59008752091509

$ws.Range("B25").Value = 1148.3014075319854
$ws.Range("D25").Value = 100.3226818908293

$ws.Range("B26").Value = 2548.1511356917727
$ws.Range("C26").Value = 0.15805396446580169
$ws.Range("D26").Value = 19.400059308747455

$ws.Range("B27").Value = 5676.1687909014872
$ws.Range("C27").Value = 0.31447368185868918
$ws.Range("D27").Value = 7.4677769386905997

$ws.Range("B28").Value = 5203.1399961988755
$ws.Range("C28").Value = 0.29778788571515025
$ws.Range("D28").Value = 5.8776189863436699
